# Generate Report for Handback
#
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   everywhere it appears (Overview sheet + each language sheet).
# - Each language sheet (zh-cn, de-de) gets its "Latest Target File" (F) and
#   "Latest Handback File" (G) columns populated for both data rows, mirroring
#   the existing "Source File Name" (A) / "Latest Handoff File" (D) hyperlinks.
# - The "Latest Handback DateTime" (H) placeholder is replaced with the real
#   handback timestamp (different per language sheet).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

function Get-HyperlinkForRange($ws, $range) {
    $target = $range.Address()
    $found = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $target) {
            $found = $hl
        }
    }
    return $found
}

function Copy-Hyperlink($ws, $fromAddr, $toAddr) {
    $source = Get-HyperlinkForRange $ws $ws.Range($fromAddr)
    if ($source -ne $null) {
        $ws.Hyperlinks.Add($ws.Range($toAddr), $source.Address, "", "", $source.TextToDisplay)
    }
}

# --- Overview sheet: refresh the status text shown per file/language ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in 2,3) {
    foreach ($col in "B","C") {
        $addr = "$col$r"
        if ($overview.Range($addr).Value() -eq $oldStatus) {
            $overview.Range($addr).Value = $newStatus
        }
    }
}

# --- Per-language detail sheets ---
$languages = @(
    @{ Sheet = "zh-cn"; HandbackDateTime = "2016-03-11 22:32:05" },
    @{ Sheet = "de-de"; HandbackDateTime = "2016-03-11 22:32:11" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    foreach ($r in 2,3) {
        # Status column
        $cAddr = "C$r"
        if ($ws.Range($cAddr).Value() -eq $oldStatus) {
            $ws.Range($cAddr).Value = $newStatus
        }

        # Source hyperlink (A) -> mirror into Latest Target File (F)
        Copy-Hyperlink $ws "A$r" "F$r"

        # Handoff-file hyperlink (D) -> mirror into Latest Handback File (G)
        Copy-Hyperlink $ws "D$r" "G$r"

        # Latest Handback DateTime
        $ws.Range("H$r").Value = $lang.HandbackDateTime
    }
}
